# logic file final update
# Update the commission/tier threshold values in column E across the
# four Sale1/Sale2 tier blocks (rows 3-6, 8-11, 13-16, 18-21, 23-26).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 3
$ws.Range("E4").Value = 4
$ws.Range("E5").Value = 5
$ws.Range("E6").Value = 6.5

$ws.Range("E8").Value = 3
$ws.Range("E9").Value = 4
$ws.Range("E10").Value = 5
$ws.Range("E11").Value = 6

$ws.Range("E13").Value = 3
$ws.Range("E14").Value = 4
$ws.Range("E15").Value = 5
$ws.Range("E16").Value = 6

$ws.Range("E18").Value = 3
$ws.Range("E19").Value = 4
$ws.Range("E20").Value = 5
$ws.Range("E21").Value = 6

$ws.Range("E23").Value = 3
$ws.Range("E24").Value = 4
$ws.Range("E25").Value = 5
$ws.Range("E26").Value = 6

# Reflect where the author's cursor ended up when the file was saved.
[void]$ws.Range("J17").Select()
